# Apply updated crypto price/volume figures to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '54.422.20'
$ws.Range("E2").Value = '  +0.33%  '
$ws.Range("D3").Value = '2.285.74'
$ws.Range("E3").Value = '  +0.10%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '502.06'
$ws.Range("E5").Value = '  +1.67%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '129.95'
$ws.Range("E6").Value = '  +2.13%  '
$ws.Range("E7").Value = '  -0.22%  '
$ws.Range("E8").Value = '  +0.42%  '
$ws.Range("E9").Value = '  +2.07%  '
$ws.Range("E10").Value = '  +0.92%  '
$ws.Range("E11").Value = '  +4.92%  '
$ws.Range("E12").Value = '  +2.36%  '
$ws.Range("D13").Value = '2.695.05'
$ws.Range("E13").Value = '  +0.29%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '22.90'
$ws.Range("E14").Value = '  +6.55%  '
$ws.Range("D15").Value = '54.380.87'
$ws.Range("E15").Value = '  +0.58%  '
$ws.Range("E16").Value = '  +0.61%  '
$ws.Range("D17").Value = '2.282.39'
$ws.Range("E17").Value = '  -0.41%  '
$ws.Range("E18").Value = '  +4.02%  '
$ws.Range("E19").Value = '  +2.95%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '304.59'
$ws.Range("E20").Value = '  +0.92%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.40'
$ws.Range("E21").Value = '  -0.05%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("E22").Value = '  -0.01%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '61.92'
$ws.Range("E23").Value = '  -2.82%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.999'
$ws.Range("E24").Value = '  -0.15%  '
$ws.Range("E25").Value = '  +2.25%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.34'
$ws.Range("E26").Value = '  +3.29%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '172.99'
$ws.Range("E27").Value = '  +4.76%  '
$ws.Range("E28").Value = '  +2.13%  '
$ws.Range("D29").Value = '0.0₃0693'
$ws.Range("E29").Value = '  +1.85%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.98'
$ws.Range("E30").Value = '  +1.81%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.09'
$ws.Range("E31").Value = '  +2.18%  '
$ws.Range("E32").Value = '  -0.07%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '17.88'
$ws.Range("E33").Value = '  +1.81%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.964'
$ws.Range("E34").Value = '  +10.94%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.997'
$ws.Range("E35").Value = '  -0.15%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.20'
$ws.Range("E36").Value = '  +1.93%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.78'
$ws.Range("E37").Value = '  +4.78%  '
$ws.Range("E38").Value = '  +0.76%  '
$ws.Range("E39").Value = '  +2.00%  '
$ws.Range("E40").Value = '  +1.54%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.99'
$ws.Range("E41").Value = '  +4.47%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '126.11'
$ws.Range("E42").Value = '  +0.16%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0497'
$ws.Range("E43").Value = '  +3.91%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0896'
$ws.Range("E44").Value = '  +1.07%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.550'
$ws.Range("E45").Value = '  +1.25%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '243.04'
$ws.Range("E46").Value = '  +2.69%  '
$ws.Range("E47").Value = '  +0.54%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0206'
$ws.Range("E48").Value = '  +1.79%  '
$ws.Range("E49").Value = '  +0.89%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '16.49'
$ws.Range("E50").Value = '  +1.74%  '
$ws.Range("E51").Value = '  -0.24%  '
